$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date column header (column AB, 28th column)
$ws.Cells.Item(1, 28).Value = "11_05_2021"

# New data values for column AB, rows 2-11
$values = @(12, 16, 34, 48, 117, 274, 406, 566, 177, 13)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 28).Value = $values[$i]
}

# The previous last-column total (AA12) used a shared formula; rewrite it
# as a standalone formula now that AB12 takes over as the shared-range end.
$ws.Range("AA12").Formula = "=SUM(AA2:AA11)"

# Totals row formula for new column
$ws.Range("AB12").Formula = "=SUM(AB2:AB11)"

# Update view: scroll right one column (old topLeftCell Q1 -> R1) and select AB12
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 18
$ws.Range("AB12").Select()
